# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# The workbook has a single sheet with a header row (row 1) and data rows 2-29.
# Column G holds the "K" stat; the underlying data pipeline recomputed these
# values (regen std/mean, calc and write s_vals) and this script writes the
# newly computed values into the existing cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 2
    6  = 3
    7  = 7
    8  = 5
    9  = 5
    10 = 3
    11 = 4
    12 = 6
    13 = 0
    14 = 4
    15 = 5
    16 = 1
    17 = 4
    18 = 8
    19 = 1
    20 = 2
    21 = 3
    22 = 5
    23 = 5
    24 = 6
    25 = 4
    26 = 5
    27 = 5
    28 = 1
    29 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
